$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "304.86"
Set-TextValue "E2" "1.01%"
Set-TextValue "D3" "35.96"
Set-TextValue "E3" "-3.98%"
Set-TextValue "D4" "5.120"
Set-TextValue "E4" "2.40%"
Set-TextValue "D5" "0.07854"
Set-TextValue "E5" "0.00%"
Set-TextValue "D6" "2.158"
Set-TextValue "E6" "-3.28%"
Set-TextValue "D7" "7.949"
Set-TextValue "E7" "-1.05%"
Set-TextValue "D8" "4.107"
Set-TextValue "E8" "2.19%"
Set-TextValue "D9" "0.9195"
Set-TextValue "E9" "1.17%"
Set-TextValue "D10" "0.09711"
Set-TextValue "E10" "2.76%"
Set-TextValue "D11" "0.1873"
Set-TextValue "E11" "-0.62%"
Set-TextValue "D12" "0.08648"
Set-TextValue "E12" "1.78%"
Set-TextValue "D13" "0.03469"
Set-TextValue "E13" "-1.52%"
Set-TextValue "D14" "0.09946"
Set-TextValue "E14" "-0.18%"
Set-TextValue "D15" "0.001443"
Set-TextValue "E15" "-2.56%"
Set-TextValue "E16" "-0.18%"
Set-TextValue "D17" "3.464"
Set-TextValue "E17" "-0.08%"
Set-TextValue "D18" "2.393"
Set-TextValue "E18" "15.17%"
Set-TextValue "D19" "0.3431"
Set-TextValue "E19" "-0.95%"
Set-TextValue "D20" "0.1319"
Set-TextValue "E20" "0.85%"
Set-TextValue "D21" "4.826"
Set-TextValue "E21" "1.36%"
Set-TextValue "D22" "0.2202"
Set-TextValue "E22" "0.00%"
Set-TextValue "D23" "0.04534"
Set-TextValue "E23" "-2.44%"
Set-TextValue "D24" "0.005084"
Set-TextValue "E24" "14.15%"
Set-TextValue "E25" "0.51%"
Set-TextValue "D26" "0.0001401"
Set-TextValue "E26" "7.81%"
Set-TextValue "D27" "0.0004753"
Set-TextValue "E27" "0.15%"
Set-TextValue "D39" "0.01848"
Set-TextValue "E39" "4.36%"
Set-TextValue "D40" "0.04780"
Set-TextValue "E40" "0.50%"
Set-TextValue "D41" "0.007801"
Set-TextValue "E41" "-0.65%"
Set-TextValue "D42" "0.1400"
Set-TextValue "E42" "0.64%"
Set-TextValue "D43" "0.007735"
Set-TextValue "E43" "0.98%"
Set-TextValue "D44" "0.002231"
Set-TextValue "E44" "0.09%"
Set-TextValue "D45" "0.01126"
Set-TextValue "E45" "14.70%"
Set-TextValue "D46" "0.00006428"
Set-TextValue "E46" "6.42%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.11%"
Set-TextValue "D48" "0.0005805"
Set-TextValue "E48" "0.07%"
Set-TextValue "D49" "47.90"
Set-TextValue "E49" "452.45%"
Set-TextValue "D50" "0.002001"
Set-TextValue "E50" "-25.54%"
Set-TextValue "D51" "0.00002101"
Set-TextValue "E51" "0.11%"
